$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.391"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05935"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.438"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8051"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9068"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1412"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07412"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03246"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03036"
$ws.Range("D13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09316"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.872"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("D16").Value = "'0.001585"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04777"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005946"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006194"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.004396"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0009865"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.00007803"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'3.608"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.148"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.03867"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006120"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1063"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002581"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007255"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005195"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.0005806"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.9591"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.002273"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
